$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 currently holds the merged "25 - Constantine" entry; split it into
# two rows: row 23 keeps the original record but with a corrected wilaya
# name, and a new row 24 captures the second record that was previously
# concatenated into C23.

$ws.Range("C23").Value = "Constantine"

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = 23
$ws.Range("C24").Value = "25 - Constantine"
$ws.Range("D24").Value = "El Khroub"
$ws.Range("E24").Value = "66670db52263a5b994b21f58"
